$wb = $excel.ActiveWorkbook

# --- Sheet: Combined_Total (rows 64-72, tax year 2023) ---
$ws = $wb.Worksheets.Item("Combined_Total")
$ws.Cells.Item(64, 4).Value = 21683.9361904972
$ws.Cells.Item(64, 6).Value = 21683.9361904972

$ws.Cells.Item(65, 4).Value = 125572.693809503
$ws.Cells.Item(65, 6).Value = 125572.693809503

$ws.Cells.Item(66, 4).Value = 19068.7296387063
$ws.Cells.Item(66, 5).Value = 13590.8188469453
$ws.Cells.Item(66, 6).Value = 32659.5484856516

$ws.Cells.Item(67, 4).Value = 79753.4903612937
$ws.Cells.Item(67, 5).Value = 56842.5511530548
$ws.Cells.Item(67, 6).Value = 136596.041514349

$ws.Cells.Item(68, 4).Value = 105386.910561649
$ws.Cells.Item(68, 5).Value = 402.716224242701
$ws.Cells.Item(68, 6).Value = 105789.626785892

$ws.Cells.Item(69, 4).Value = 747650.279629641
$ws.Cells.Item(69, 5).Value = 2728.61384526203
$ws.Cells.Item(69, 6).Value = 750378.893474903

$ws.Cells.Item(70, 4).Value = 196577.75980871
$ws.Cells.Item(70, 5).Value = 746.939930495271
$ws.Cells.Item(70, 6).Value = 197324.699739205

$ws.Cells.Item(71, 4).Value = 13227.6073808706
$ws.Cells.Item(71, 6).Value = 13227.6073808706

$ws.Cells.Item(72, 4).Value = 127121.182619129
$ws.Cells.Item(72, 6).Value = 127121.182619129

# --- Sheet: Combined_Columbus (rows 20-21, tax year 2023) ---
$ws = $wb.Worksheets.Item("Combined_Columbus")
$ws.Cells.Item(20, 4).Value = 21683.9361904972
$ws.Cells.Item(20, 6).Value = 21683.9361904972

$ws.Cells.Item(21, 4).Value = 125572.693809503
$ws.Cells.Item(21, 6).Value = 125572.693809503

# --- Sheet: Combined_Gahanna (rows 20-21, tax year 2023) ---
$ws = $wb.Worksheets.Item("Combined_Gahanna")
$ws.Cells.Item(20, 4).Value = 19068.7296387063
$ws.Cells.Item(20, 5).Value = 13590.8188469453
$ws.Cells.Item(20, 6).Value = 32659.5484856516

$ws.Cells.Item(21, 4).Value = 79753.4903612937
$ws.Cells.Item(21, 5).Value = 56842.5511530548
$ws.Cells.Item(21, 6).Value = 136596.041514349

# --- Sheet: Combined_JeffersonUnincorporate (rows 20-22, tax year 2023) ---
$ws = $wb.Worksheets.Item("Combined_JeffersonUnincorporate")
$ws.Cells.Item(20, 4).Value = 105386.910561649
$ws.Cells.Item(20, 5).Value = 402.716224242701
$ws.Cells.Item(20, 6).Value = 105789.626785892

$ws.Cells.Item(21, 4).Value = 747650.279629641
$ws.Cells.Item(21, 5).Value = 2728.61384526203
$ws.Cells.Item(21, 6).Value = 750378.893474903

$ws.Cells.Item(22, 4).Value = 196577.75980871
$ws.Cells.Item(22, 5).Value = 746.939930495271
$ws.Cells.Item(22, 6).Value = 197324.699739205

# --- Sheet: Combined_Reynoldsburg (rows 10-11, tax year 2023) ---
$ws = $wb.Worksheets.Item("Combined_Reynoldsburg")
$ws.Cells.Item(10, 4).Value = 13227.6073808706
$ws.Cells.Item(10, 6).Value = 13227.6073808706

$ws.Cells.Item(11, 4).Value = 127121.182619129
$ws.Cells.Item(11, 6).Value = 127121.182619129
